$wb = $excel.ActiveWorkbook

# InvoicePage sheet gets two new rows of data
$ws = $wb.Worksheets.Item("InvoicePage")

$ws.Range("A4").Value = "ItemQuantity"
$ws.Range("B4").Value = 5

$ws.Range("A5").Value = "Rate"
$ws.Range("B5").Value = 55

$ws.Columns.Item(1).ColumnWidth = 12.44140625

$ws.Activate()
$ws.Range("B5").Select()

$wb.Save()
